$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Insert a new column in front of the "50%_germ" column (AD) and give it
# the new header "germ_rate_days" -- this shifts the old AD ("50%_germ")
# to AE and the old AE ("Notes") to AF, matching the source edit.
$ws.Columns("AD:AD").Insert()
$ws.Cells.Item(1, 30).Value = "germ_rate_days"

# Match the row height Excel recalculated for the (now taller) header row.
$ws.Rows(1).RowHeight = 70

# Restore the view/selection state recorded in the saved workbook: the
# user had scrolled right and landed just past the new data.
$ws.Range("R1").Select() | Out-Null
$ws.Range("AG1").Select() | Out-Null
